$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename several sheet tabs (commit message: "change several spreadsheet
#    tab names, and one column title")
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("total mortality").Name        = "mortality rates"
$wb.Worksheets.Item("mortality").Name               = "causes of death"
$wb.Worksheets.Item("RRStunting").Name              = "RR death by stunting"
$wb.Worksheets.Item("RRWasting").Name               = "RR death by wasting"
$wb.Worksheets.Item("RRBreastfeeding").Name         = "RR death by breastfeeding"
$wb.Worksheets.Item("RR Death by Birth Outcome").Name = "RR death by birth outcome"
$wb.Worksheets.Item("OR stunting for complements").Name = "OR stunting by compfeeding"
$wb.Worksheets.Item("OR appropriateBF by interv").Name  = "OR correctBF by interventn"

# ---------------------------------------------------------------------------
# 2. Change the column title in "OR stunting by compfeeding" (formerly
#    "OR stunting for complements") from "Complements group" to
#    "Food security & education".
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("OR stunting by compfeeding")
$wsComp.Range("A1").Value = "Food security & education"

# ---------------------------------------------------------------------------
# 3. Update the active sheet / selection so that "OR stunting by compfeeding"
#    becomes the selected tab with A9 selected, matching the new view state.
# ---------------------------------------------------------------------------
$wsComp.Activate()
$wsComp.Range("A9").Select()
